$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Codelists_2MD")

# Delete the "BSTED" codelist rows (originally rows 122-124) first,
# then the "bosted" codelist rows (originally rows 102-107), bottom-up
# so row numbers don't shift under us.
$ws.Rows("122:124").Delete()
$ws.Rows("102:107").Delete()

# Update view state to match target: Codelists_2MD becomes the active/selected sheet.
$ws.Activate()
$ws.Application.ActiveWindow.ScrollRow = 82
$ws.Range("A117").Select()

$ws1 = $wb.Worksheets.Item("Variables_MD")
$ws1.Range("E25").Select()
